$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Linear)
$ws.Range("B2").Value = 1.955880088503596
$ws.Range("C2").Value = 3.359480197017807
$ws.Range("D2").Value = 0.7524606424109233

# Row 3 (Decision Tree)
$ws.Range("B3").Value = 1.977707360461119
$ws.Range("C3").Value = 3.145172413793103
$ws.Range("D3").Value = 0.7412243088690889

# Row 4 (Random Forest)
$ws.Range("B4").Value = 1.487185414381639
$ws.Range("C4").Value = 1.709907722195245
$ws.Range("D4").Value = 0.9172561303656908

# Row 5 (Lasso)
$ws.Range("B5").Value = 1.952506025388216
$ws.Range("C5").Value = 3.338063459239361
$ws.Range("D5").Value = 0.7541643351512322
